$wb = $excel.ActiveWorkbook

# --- Step 1: update the view state of "13th January 2022-v2" (it is no longer the active/last sheet) ---
$ws5 = $wb.Worksheets.Item("13th January 2022-v2")
$ws5.Activate()
$ws5.Range("A1:A115").Select()
$excel.ActiveWindow.ScrollRow = 85

# --- Step 2: add the new sheet "16th January 2022-v1" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws6.Name = "16th January 2022-v1"

# --- Step 3: populate column A of the new sheet (row -> text), mirroring the mod load order list ---
$data = @(
    @(1, 'STALKER Anomaly 1.5.1'),
    @(2, '16th January 2022'),
    @(3, 'Mod "Load Order" version 1'),
    @(5, 'Tools - MCM 1.5'),
    @(6, 'Tools - MCM 1.6'),
    @(7, 'Performance - Alife Optimized'),
    @(8, 'Performance - Crowded areas'),
    @(9, 'Gamepaly - 00. Main EFT reposition'),
    @(10, 'Weapons Sounds - JSRS P.1 - Main Addon'),
    @(11, 'Weapons Sounds - JSRS P.2 - Cracks + Impact Patch'),
    @(12, 'Weapons Sounds - JSRS P.3 - New Tracers Patch'),
    @(13, 'QOL - AdjustableScope 1.6'),
    @(14, 'Weapons Sounds - Fire Mode Selector Sound'),
    @(15, 'Weapons Sounds - Fire Mode Selector Sound - JSRS'),
    @(16, 'Animation - I.N.E.R.T.I.A'),
    @(17, 'Animation - Basic Animation Reworks'),
    @(18, 'Animation - Lower Weapon Sprint Reworks'),
    @(19, 'Animation - Death Animation Rework'),
    @(20, 'Animation - Food_Drug_Drinks Animate'),
    @(21, 'BG Sounds - SCO3 P.1 Real Mutant Sounds'),
    @(22, 'BG Sounds - SCO3 P.2 New Outdoor Sounds'),
    @(23, 'BG Sounds - SCO3 P.3 New Underground Sounds'),
    @(24, 'BG Sounds - SCO3 P.4 Dead Zone (No Animal Sounds)'),
    @(25, 'Misc Sounds-EFTfootstep1.1'),
    @(26, 'Misc Sounds-EFTJumpLandSFX_HarukaSai'),
    @(27, 'Misc Sounds-EFTAimRattle1.2'),
    @(28, 'Misc Sounds-Inventory Open and Close Sound'),
    @(29, 'Misc Sounds-Quiter Wood Breaking'),
    @(30, 'Misc Sounds-Exo Servomotor Sound'),
    @(31, 'Misc Sounds-Chimaera Growls'),
    @(32, 'Particle Effects - OPTION 2 - MinimalisticFX'),
    @(33, 'Particle Effects - 00.CinematicVFX3.5UPD1.1 - Core'),
    @(34, 'Particle Effects - 01.CVFX Blood Effects Tweak - Stronger'),
    @(35, 'Textures - Vanilla Models Re_Textures 3.7'),
    @(36, 'Textures - UNISG Overhaul'),
    @(37, 'Textures - SEVA Glass Variety'),
    @(38, 'Textures - (Optional) Grass Yellow Reeds'),
    @(39, 'Gameplay - Better NewGame Loadouts'),
    @(40, 'Gameplay - Trader Overhaul (Main Folder)'),
    @(41, 'Gameplay - RPGlite Skills Systems 1.4.2.3'),
    @(42, 'Gameplay - Artefact Reworks 1.0.4'),
    @(43, 'Gameplay - 00. Grok''s Stash Overhaul'),
    @(44, 'Gameplay - 01. TB Fix Bugged Stashes Patch'),
    @(45, 'Gameplay - GBOOBS Core'),
    @(46, 'Gameplay - GBOOBS Affects all NPCs'),
    @(47, 'Gameplay - GBOOBS Shotguns fix'),
    @(48, 'Gameplay - AI use cover Beta v.3.1'),
    @(49, 'Gameplay - Ignore Companion Allegiance'),
    @(50, 'Gameplay - Knife in inventory to skin mutants'),
    @(51, 'HUD - 00 SIDHUD MAIN FILE [INSTALL FIRST ALWAYS]'),
    @(52, 'HUD - Cr3pis Icon Rework Vanilla Models (with labels)'),
    @(53, 'HUD - Instant Tooltip'),
    @(54, 'HUD - Hunger Thirst Sleep UI 0.61'),
    @(55, 'HUD - BodyDotsMiniMap'),
    @(56, 'HUD - BatteryWarning'),
    @(57, 'HUD - AscendantCrosshairs'),
    @(58, 'HUD - Tacticool Scopes DX8'),
    @(59, 'Quests - LTTZ DC - 1.5.1'),
    @(60, 'Quests - LTTZ DC - MotZ Unlocker'),
    @(61, 'Quests - MoTZ dialogue unlocker'),
    @(62, 'QOL - SortingItems'),
    @(63, 'QOL - Stealth_2.0'),
    @(64, 'QOL - Persistent Weather (Yohji) 1.3'),
    @(65, 'QOL - No Grey RadiationPostProcess'),
    @(66, 'QOL - Arszis Mutant Bleeding 2.0.1'),
    @(67, 'QOL - Remove Mutant AutoAim'),
    @(68, 'QOL - ARX and ABF Anomalous Rework'),
    @(69, 'QOL - BasePropagandaReturn'),
    @(70, 'QOL - Bottle Water Adjusted'),
    @(71, 'QOL - Coordinate Based Safe Zone'),
    @(72, 'QOL - CozyCampfires'),
    @(73, 'QOL - Dialogue Expanded 4.1'),
    @(74, 'QOL - Enable Combine Items of Similar Kind into One'),
    @(75, 'QOL - Enable Disassemble All Items'),
    @(76, 'QOL - Encyclopedia_ui_fixes_v1.2'),
    @(77, 'QOL - Food Drugs Drinks Visual Accuracy Item Counts'),
    @(78, 'QOL - Highlight New Items'),
    @(79, 'QOL - Keep Crafting Window Open'),
    @(80, 'QOL - NVR alternative'),
    @(81, 'QOL - NVR EFT Style Version Old DX10-9-8'),
    @(82, 'QOL - NVR EFT Style Version Old DX10-9-8 (Clean-ZM Modification)'),
    @(83, 'Misc Sounds - EFTNOD'),
    @(84, 'HUD - Collimator Sights Rework'),
    @(85, 'HUD - Subtle Hit Direction Marker 1.1'),
    @(86, 'QOL - NicerFlashlights Addon Taclight 1.6'),
    @(87, 'QOL - Disable Goodwill Loss'),
    @(88, 'QOL - NPC''s Loadout Rework - 00. Main'),
    @(89, 'QOL - NPC''s Loadout Rework - 02. Advanced Army'),
    @(90, 'QOL - NPC''s Loadout Rework - 02. Improved Army'),
    @(91, 'QOL - NPC''s Loadout Rework - 04. Advanced Monolith'),
    @(92, 'QOL - NPC''s Loadout Rework - 05. Advanced Merc'),
    @(93, 'QOL - NPC''s Loadout Rework - 06. Advanced Ecologists'),
    @(94, 'QOL - 00 Dynamic_NPC_Armor_Visuals'),
    @(95, 'QOL - 01 Dynamic_NPC_Armor_Visuals NPCs_loot_outfits'),
    @(96, 'QOL - Bounty Squad Rework'),
    @(97, 'QOL - Better Companion HP Regen (10x)'),
    @(98, 'QOL - CompanionInventoryUnlock 04_01_2022'),
    @(99, 'QOL - Companions Deactive Headlamp'),
    @(100, 'QOL - CompanionsDon''tDie'),
    @(101, 'QOL - Companion Carry Weights Tweak'),
    @(102, 'QOL - Decreased AI Hearing Distance (Bushes)'),
    @(103, 'QOL - European UNISG Names v3'),
    @(104, 'QOL - Merc Codenames LTTZ DC compatible'),
    @(105, 'QOL - Quest Items Droppable RC 18'),
    @(106, 'QOL - Suppressor Reworked 50 - Standard Anomaly 1.05.6'),
    @(107, 'QOL - TacticoolScopes_RadiusFix'),
    @(108, 'QOL - ToggleScope'),
    @(109, 'QOL - Weapon Sway 0.5'),
    @(110, 'QOL - Less NPC Jamming'),
    @(111, 'QOL - Fluid_Aim_v1.2.3'),
    @(112, 'QOL - WeaponPartsRework - 2022_01_08'),
    @(113, 'QOL - Upgraded Weapons Weight Rework - 2021_12_31'),
    @(114, 'QOL - Unjam Keybind - 2022_01_08'),
    @(115, 'Armor - Sarcophagus'),
    @(116, 'QOL - Exo Missing Texture Fix'),
    @(117, 'QOL - Outfit Ballistic Upgrade Rework'),
    @(118, 'QOL - Outfit Speed Rework'),
    @(119, 'QOL - Reworked Outfit Attachments 1.3'),
    @(120, 'HUD - OPT1. Groks New Masks Reflections Droplets 2.1.0'),
    @(121, 'Weapons - XM8R - Main Addon'),
    @(122, 'Weapons - XM8R - patch JSRS 3.6'),
    @(123, 'Weapons Sounds - Suppressed Shotgun Sound Rework JSRS'),
    @(124, 'Gameplay - Blindsides RPL 0.7'),
    @(125, 'Gameplay - 02. EFT repos + Blindsides reanim (vanilla stats)'),
    @(126, 'Gameplay - Blindsides RPL GBOOBS Patch'),
    @(127, 'Gameplay - AWAR 0.97'),
    @(128, 'Gameplay - 02. EFT repos + Blindsides reanim'),
    @(129, 'Gameplay - 01. EFT repos + AWAR 0.97'),
    @(130, 'Gameplay - Blindsides RPL JSRS Patch 0.6')
)

foreach ($pair in $data) {
    $r = [int]$pair[0]
    $v = [string]$pair[1]
    $ws6.Cells.Item($r, 1).Value = $v
}

# --- Step 4: set the view state of the new sheet (it becomes the active/selected tab) ---
$ws6.Activate()
$ws6.Range("A131").Select()
$excel.ActiveWindow.ScrollRow = 104

Write-Output "Added sheet '16th January 2022-v1' with $($data.Count) populated rows"
